$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("利欧股份","金风科技","金风科技"),
    @("蓝色光标","岩山科技","海格通信"),
    @("岩山科技","利欧股份","岩山科技"),
    @("特变电工","中国卫星","利欧股份"),
    @("海格通信","蓝色光标","特变电工"),
    @("三维通信","航天发展","航天发展"),
    @("航天发展","航天电子","鲁信创投"),
    @("金风科技","三维通信","蓝色光标"),
    @("鲁信创投","雷科防务","中国卫星"),
    @("中国卫星","特变电工","通宇通讯"),
    @("航天电子","海格通信","航天电子"),
    @("美年健康","山子高科","雷科防务"),
    @("雷科防务","中国卫通","三维通信"),
    @("通宇通讯","通宇通讯","巨力索具"),
    @("易点天下","易点天下","省广集团"),
    @("山子高科","鲁信创投","东百集团"),
    @("中国卫通","美年健康","吉视传媒"),
    @("巨力索具","省广集团","雪人集团"),
    @("雪人集团","雪人集团","平潭发展"),
    @("省广集团","巨力索具","山子高科")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
